# Applies the reordering of Category (A) / NumberOfPrisoners (B) values
# for specific rows in the Prisoners_Data worksheet, as described by the
# target diff. Column C (Year) and D (Country) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = @{ A = "F"; B = 16 }
    10 = @{ A = "M"; B = 5033 }
    12 = @{ A = "C"; B = 156 }
    13 = @{ A = "M"; B = 6200 }
    20 = @{ A = "C"; B = 350 }
    21 = @{ A = "M"; B = 6171 }
    22 = @{ A = "F"; B = 58 }
    23 = @{ A = "F"; B = 54 }
    24 = @{ A = "M"; B = 5500 }
    25 = @{ A = "C"; B = 230 }
    26 = @{ A = "C"; B = 185 }
    27 = @{ A = "M"; B = 5000 }
    29 = @{ A = "M"; B = 4400 }
    30 = @{ A = "C"; B = 170 }
    31 = @{ A = "F"; B = 41 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
}
